# Generate Report for Handback
# Update the "Latest Handback DateTime" (column K, row 2) on the
# per-language handback-status sheets to reflect the newly received
# handback files for fd49981e-6c8d-4974-addf-06909b19656f.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-11-09 01:20:13"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-11-09 01:20:33"
